$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New label text for column B (rows 2-31)
$newLabel = "Índice de volume de vendas no comércio varejista - Variação mensal (base: igual mês do ano anterior)"

# New values for column D (rows 2-31)
$dValues = @{
    2  = 5.1
    3  = 3.9
    4  = 0.3
    5  = -7.2
    6  = -4.9
    7  = 4
    8  = 0.7
    9  = 2.6
    10 = 1.3
    11 = -2.9
    12 = 5.6
    13 = 5.122222222222223
    14 = 0.4111111111111111
    15 = -10.13333333333333
    16 = -5.233333333333333
    17 = 1.077777777777778
    18 = -0.4555555555555555
    19 = 1.311111111111111
    20 = 2.855555555555556
    21 = -8.377777777777778
    22 = 2.3
    23 = -0.2
    24 = 0.5
    25 = -10.6
    26 = -1.4
    27 = -9.5
    28 = 0.6
    29 = -1.8
    30 = 0.2
    31 = -11.2
}

for ($row = 2; $row -le 31; $row++) {
    $ws.Cells.Item($row, 2).Value = $newLabel
    $ws.Cells.Item($row, 4).Value = $dValues[$row]
}
